# Update Design and Plan documents
# Applies the edits captured in the target diff:
#  - Row height changes for the risk table rows (2-6), caused by the new,
#    longer wrapped text rendering on a different machine/font metrics.
#  - Two new "Task" rows (10 and 11) filled out with owner/start/due/status.
#  - A new date-formatted cell style (center/wrap/shrink, d-mmm format)
#    used by the new Start Date / Due Date cells.
#  - Updated selection/active cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights for rows 2-6 (auto-fit-like re-wrap heights from the diff) ---
$ws.Rows.Item(2).RowHeight = 111.75
$ws.Rows.Item(3).RowHeight = 121.5
$ws.Rows.Item(4).RowHeight = 99.75
$ws.Rows.Item(5).RowHeight = 111.75
$ws.Rows.Item(6).RowHeight = 99

# --- Row 10: new task "Task Title" (A10 already has text) gets Owner/Start/Due/Complete ---
$ws.Range("B10").Value = "Jayden Brooks"
$ws.Range("C10").Value = 45991
$ws.Range("D10").Value = 46011
$ws.Range("E10").Value = "Complete"

# --- Row 11: same pattern ---
$ws.Range("B11").Value = "Jayden Brooks"
$ws.Range("C11").Value = 45991
$ws.Range("D11").Value = 46011
$ws.Range("E11").Value = "Complete"

# --- Date formatting/style for the new Start Date / Due Date cells ---
$dateRange = $ws.Range("C10:D11")
$dateRange.NumberFormat = "d-mmm"
$dateRange.HorizontalAlignment = -4108
$dateRange.VerticalAlignment = -4108
$dateRange.WrapText = $true
$dateRange.ShrinkToFit = $true

# --- Update the sheet's active selection, as recorded in the saved view state ---
$ws.Range("G15").Select()
